$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("blueprint_tasks")

# Fix casing of two text values (shared strings used by these cells)
$ws.Range("B8").Value = "roomplan-FP"
$ws.Range("C22").Value = "Walkthrough-video"

# Update the view: scroll so column C is the leftmost visible column,
# and move the active cell/selection to C22 (was B8)
try {
    $excel.ActiveWindow.ScrollColumn = 3
} catch {
}
$ws.Range("C22").Select()
